# Update loading_percent values for case with 380 kV (Case_2_32)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 11.99252434514551
$ws.Range("C2").Value = 8.253053475413243
$ws.Range("E2").Value = 16.31847163664813
$ws.Range("F2").Value = 16.86991607391233
$ws.Range("G2").Value = 3.671679663212849
$ws.Range("I2").Value = 25.45930078438566
$ws.Range("K2").Value = 12.56597086443528
$ws.Range("N2").Value = 20.10175556219518

# Row 3
$ws.Range("B3").Value = 11.71418294256921
$ws.Range("C3").Value = 7.989239252821621
$ws.Range("E3").Value = 15.40303801569763
$ws.Range("F3").Value = 15.89584955866808
$ws.Range("G3").Value = 3.674756691367645
$ws.Range("I3").Value = 25.47667431675478
$ws.Range("K3").Value = 12.36731976164734
$ws.Range("N3").Value = 20.1616107516428

# Row 4
$ws.Range("B4").Value = 11.54364770689598
$ws.Range("C4").Value = 7.825750856553686
$ws.Range("E4").Value = 14.81763563796413
$ws.Range("F4").Value = 15.26997757108491
$ws.Range("G4").Value = 3.67674073998142
$ws.Range("I4").Value = 25.49289254485977
$ws.Range("K4").Value = 12.24754701713519
$ws.Range("N4").Value = 20.20031657301585

# Row 5
$ws.Range("B5").Value = 11.47435791801514
$ws.Range("C5").Value = 7.758867006367395
$ws.Range("E5").Value = 14.57348992066987
$ws.Range("F5").Value = 15.00819731993403
$ws.Range("G5").Value = 3.677573173824426
$ws.Range("I5").Value = 25.50089197829189
$ws.Range("K5").Value = 12.19935699253761
$ws.Range("N5").Value = 20.21658124861571

# Row 6
$ws.Range("B6").Value = 11.46286809852263
$ws.Range("C6").Value = 7.747748807506741
$ws.Range("E6").Value = 14.53262119973611
$ws.Range("F6").Value = 14.96433081551593
$ws.Range("G6").Value = 3.677712846187589
$ws.Range("I6").Value = 25.50230408415776
$ws.Range("K6").Value = 12.19139430648562
$ws.Range("N6").Value = 20.2193116936986

# Row 7
$ws.Range("B7").Value = 11.54271226143597
$ws.Range("C7").Value = 7.824849733675326
$ws.Range("E7").Value = 14.8143652326184
$ws.Range("F7").Value = 15.26647399323137
$ws.Range("G7").Value = 3.676751869509356
$ws.Range("I7").Value = 25.49299480617859
$ws.Range("K7").Value = 12.24689452067156
$ws.Range("N7").Value = 20.20053393280291

# Row 8
$ws.Range("B8").Value = 11.89654240886094
$ws.Range("C8").Value = 8.16247097690902
$ws.Range("E8").Value = 16.00782013078562
$ws.Range("F8").Value = 16.53996406344768
$ws.Range("G8").Value = 3.672721016135536
$ws.Range("I8").Value = 25.46413600058789
$ws.Range("K8").Value = 12.49705875343292
$ws.Range("N8").Value = 20.12198813482953

# Row 9
$ws.Range("B9").Value = 12.5885532399939
$ws.Range("C9").Value = 8.807830572332268
$ws.Range("E9").Value = 18.18007734848568
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 3.665563925238366
$ws.Range("I9").Value = 25.45180892132626
$ws.Range("K9").Value = 13.00219010342267
$ws.Range("N9").Value = 19.9834552683365

# Row 10
$ws.Range("B10").Value = 13.08977314627468
$ws.Range("C10").Value = 9.265988471893058
$ws.Range("E10").Value = 19.80377231681776
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 3.660755172336625
$ws.Range("I10").Value = 25.47000989727561
$ws.Range("K10").Value = 13.37831398187437
$ws.Range("N10").Value = 19.89110266075042

# Row 11
$ws.Range("B11").Value = 13.31504800084631
$ws.Range("C11").Value = 9.469891166364768
$ws.Range("E11").Value = 20.50116567895176
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 3.658663862892237
$ws.Range("I11").Value = 25.48425491169316
$ws.Range("K11").Value = 13.54970424552934
$ws.Range("N11").Value = 19.85113271554965

# Row 12
$ws.Range("B12").Value = 13.3998707028657
$ws.Range("C12").Value = 9.546378194727627
$ws.Range("E12").Value = 20.7593588862059
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 3.657885674085342
$ws.Range("I12").Value = 25.49050955137065
$ws.Range("K12").Value = 13.61458193795215
$ws.Range("N12").Value = 19.83629062599436

# Row 13
$ws.Range("B13").Value = 13.38162552352128
$ws.Range("C13").Value = 9.529938797471129
$ws.Range("E13").Value = 20.70401396746519
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 3.658052661117788
$ws.Range("I13").Value = 25.48912420421284
$ws.Range("K13").Value = 13.60061145980892
$ws.Range("N13").Value = 19.8394740734234

# Row 14
$ws.Range("B14").Value = 13.32203661427886
$ws.Range("C14").Value = 9.476198773161489
$ws.Range("E14").Value = 20.52252532573725
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 3.658599565887715
$ws.Range("I14").Value = 25.48475222942676
$ws.Range("K14").Value = 13.5550426235172
$ws.Range("N14").Value = 19.84990576200277

# Row 15
$ws.Range("B15").Value = 13.28547106351761
$ws.Range("C15").Value = 9.443184695243323
$ws.Range("E15").Value = 20.41059168159317
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 3.658936348114746
$ws.Range("I15").Value = 25.48218637781687
$ws.Range("K15").Value = 13.52712534649566
$ws.Range("N15").Value = 19.85633371837093

# Row 16
$ws.Range("B16").Value = 13.07498798832531
$ws.Range("C16").Value = 9.252565301736187
$ws.Range("E16").Value = 19.75736927556254
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 3.660893772653355
$ws.Range("I16").Value = 25.46919922293176
$ws.Range("K16").Value = 13.36711338509966
$ws.Range("N16").Value = 19.89375587993014

# Row 17
$ws.Range("B17").Value = 12.94509674611191
$ws.Range("C17").Value = 9.134412393397209
$ws.Range("E17").Value = 19.34610240495273
$ws.Range("F17").Value = 20.20408069617459
$ws.Range("G17").Value = 3.662119167635795
$ws.Range("I17").Value = 25.4627619095118
$ws.Range("K17").Value = 13.26897797200123
$ws.Range("N17").Value = 19.91723613112918

# Row 18
$ws.Range("B18").Value = 12.87013457957146
$ws.Range("C18").Value = 9.066033009676126
$ws.Range("E18").Value = 19.10566940421378
$ws.Range("F18").Value = 19.95656407809808
$ws.Range("G18").Value = 3.662833044194968
$ws.Range("I18").Value = 25.45962070792021
$ws.Range("K18").Value = 13.21256414768804
$ws.Range("N18").Value = 19.93093352075319

# Row 19
$ws.Range("B19").Value = 12.84471324846822
$ws.Range("C19").Value = 9.042811111907234
$ws.Range("E19").Value = 19.02359491878548
$ws.Range("F19").Value = 19.87204792380562
$ws.Range("G19").Value = 3.663076309678285
$ws.Range("I19").Value = 25.45865347816484
$ws.Range("K19").Value = 13.19347077006789
$ws.Range("N19").Value = 19.9356042255335

# Row 20
$ws.Range("B20").Value = 12.95895066815492
$ws.Range("C20").Value = 9.147034124533494
$ws.Range("E20").Value = 19.39028410535412
$ws.Range("F20").Value = 20.24955283636157
$ws.Range("G20").Value = 3.661987785053439
$ws.Range("I20").Value = 25.46338905243978
$ws.Range("K20").Value = 13.27942192180046
$ws.Range("N20").Value = 19.91471672923823

# Row 21
$ws.Range("B21").Value = 13.33955315065696
$ws.Range("C21").Value = 9.492003812231598
$ws.Range("E21").Value = 20.57599263827937
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 3.658438554398042
$ws.Range("I21").Value = 25.48601301729293
$ws.Range("K21").Value = 13.56842846166511
$ws.Range("N21").Value = 19.84683375161137

# Row 22
$ws.Range("B22").Value = 13.58544009858129
$ws.Range("C22").Value = 9.713196763976606
$ws.Range("E22").Value = 21.31657993932416
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 3.656198999436962
$ws.Range("I22").Value = 25.50581465336873
$ws.Range("K22").Value = 13.75714173739411
$ws.Range("N22").Value = 19.80418004086239

# Row 23
$ws.Range("B23").Value = 13.45449559965041
$ws.Range("C23").Value = 9.595555536253464
$ws.Range("E23").Value = 20.9244453774756
$ws.Range("F23").Value = 21.82633154475857
$ws.Range("G23").Value = 3.657386995482006
$ws.Range("I23").Value = 25.4947865480721
$ws.Range("K23").Value = 13.65645851584139
$ws.Range("N23").Value = 19.8267884735701

# Row 24
$ws.Range("B24").Value = 12.95268819368485
$ws.Range("C24").Value = 9.141329240899132
$ws.Range("E24").Value = 19.37032199031024
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 3.66204715389134
$ws.Range("I24").Value = 25.46310377784116
$ws.Range("K24").Value = 13.27470018980359
$ws.Range("N24").Value = 19.91585513315389

# Row 25
$ws.Range("B25").Value = 12.40219062007885
$ws.Range("C25").Value = 8.635665055990643
$ws.Range("E25").Value = 17.59503704796536
$ws.Range("F25").Value = 18.34778573295691
$ws.Range("G25").Value = 3.667420722324609
$ws.Range("I25").Value = 25.45037433351958
$ws.Range("K25").Value = 12.86438211941895
$ws.Range("N25").Value = 20.01927528353927
